# Update column G ("K") values for rows 2-31 on the active sheet.
# These values were regenerated (per commit message: "regen save_data to use
# K instead of Strike#, regen std/mean, calc and write s_vals") - the net
# effect on this sheet is that column G holds new computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 2
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 4
    21 = 2
    22 = 3
    23 = 3
    24 = 2
    25 = 3
    26 = 4
    27 = 0
    28 = 1
    29 = 3
    30 = 4
    31 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
